$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 279, shifting existing rows 279-296 down to 280-297.
$ws.Rows.Item(279).Insert()

# Fill the newly inserted row 279 with the new data entry.
$ws.Cells.Item(279, 1).Value = 4
$ws.Cells.Item(279, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(279, 3).Value = "Los Lagos"
$ws.Cells.Item(279, 4).Value = 44706
$ws.Cells.Item(279, 5).Value = 10
$ws.Cells.Item(279, 6).Value = 100114014
$ws.Cells.Item(279, 7).Value = "Betarraga"
$ws.Cells.Item(279, 8).Value = "Sin especificar"
$ws.Cells.Item(279, 9).Value = "Primera"
$ws.Cells.Item(279, 10).Value = 250
$ws.Cells.Item(279, 11).Value = 1000
$ws.Cells.Item(279, 12).Value = 1000
$ws.Cells.Item(279, 13).Value = 1000
$ws.Cells.Item(279, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(279, 15).Value = "Región del Maule"
$ws.Cells.Item(279, 16).Value = 200
$ws.Cells.Item(279, 17).Value = 5
$ws.Cells.Item(279, 18).Value = "Hortaliza"
